# Fix the ordering of a couple of rows whose "id" (column B) values were
# out of sequence: rows 149/150 need to be swapped, and rows 286/287/288
# need to be rotated, so that column B (id) is sorted ascending again.
# Columns A (row index), C (Div) and D (Date) are identical / sequential
# and stay untouched; columns B and E:AD move together with each record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 149 and 150 ---
$b149 = $ws.Range("B149").Value2
$b150 = $ws.Range("B150").Value2
$rest149 = $ws.Range("E149:AD149").Value2
$rest150 = $ws.Range("E150:AD150").Value2

$ws.Range("B149").Value2 = $b150
$ws.Range("E149:AD149").Value2 = $rest150

$ws.Range("B150").Value2 = $b149
$ws.Range("E150:AD150").Value2 = $rest149

# --- Rotate rows 286, 287, 288 ---
# New 286 = old 288, New 287 = old 286, New 288 = old 287
$b286 = $ws.Range("B286").Value2
$b287 = $ws.Range("B287").Value2
$b288 = $ws.Range("B288").Value2
$rest286 = $ws.Range("E286:AD286").Value2
$rest287 = $ws.Range("E287:AD287").Value2
$rest288 = $ws.Range("E288:AD288").Value2

$ws.Range("B286").Value2 = $b288
$ws.Range("E286:AD286").Value2 = $rest288

$ws.Range("B287").Value2 = $b286
$ws.Range("E287:AD287").Value2 = $rest286

$ws.Range("B288").Value2 = $b287
$ws.Range("E288:AD288").Value2 = $rest287
